$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 8
$ws.Range("C5").Value = 17
$ws.Range("C6").Value = 11
$ws.Range("C7").Value = 23
$ws.Range("C8").Value = 17
$ws.Range("C9").Value = 12
$ws.Range("C11").Value = 14
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 13
$ws.Range("C15").Value = 16
$ws.Range("C16").Value = 19
$ws.Range("C17").Value = 18
$ws.Range("C18").Value = 19

# Update text values in column B
$ws.Range("B10").Value = "<hind>"
$ws.Range("B18").Value = "<uniform>"
